$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ============================================================================
# Step 1: copy cell FORMATTING ONLY (xlPasteFormats) from existing donor cells
# that already carry the exact target style index, onto the new cells in rows
# 100-103. Must happen BEFORE we repurpose P96 (donor for P100's style) below.
# ============================================================================
$ws1.Range("A6").Copy()
$ws1.Range("A100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("B26").Copy()
$ws1.Range("B100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("C26").Copy()
$ws1.Range("C100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("D26").Copy()
$ws1.Range("D100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("E26").Copy()
$ws1.Range("E100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("D26").Copy()
$ws1.Range("F100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("G31").Copy()
$ws1.Range("G100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("H6").Copy()
$ws1.Range("H100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("I26").Copy()
$ws1.Range("I100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("D26").Copy()
$ws1.Range("J100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("I26").Copy()
$ws1.Range("K100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("D26").Copy()
$ws1.Range("L100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("G31").Copy()
$ws1.Range("M100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("H6").Copy()
$ws1.Range("N100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P96").Copy()
$ws1.Range("P100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws1.Range("P26").Copy()
$ws1.Range("A101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("B27").Copy()
$ws1.Range("B101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("C101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("D101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("E101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("F27").Copy()
$ws1.Range("F101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("G7").Copy()
$ws1.Range("G101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("H3").Copy()
$ws1.Range("H101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("I101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("F27").Copy()
$ws1.Range("J101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("K101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("F27").Copy()
$ws1.Range("L101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("M78").Copy()
$ws1.Range("M101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("H3").Copy()
$ws1.Range("N101").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws1.Range("P26").Copy()
$ws1.Range("A102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("B27").Copy()
$ws1.Range("B102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("C102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("D102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("E102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("F102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("G102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("B27").Copy()
$ws1.Range("H102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("I102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("F27").Copy()
$ws1.Range("J102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("K102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("F27").Copy()
$ws1.Range("L102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("G7").Copy()
$ws1.Range("M102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("H3").Copy()
$ws1.Range("N102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws1.Range("P26").Copy()
$ws1.Range("A103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("B27").Copy()
$ws1.Range("B103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("C103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("D103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("E103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("F103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("G103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("B27").Copy()
$ws1.Range("H103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("I103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("F27").Copy()
$ws1.Range("J103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("K103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("F27").Copy()
$ws1.Range("L103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P26").Copy()
$ws1.Range("M103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("M23").Copy()
$ws1.Range("N103").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ============================================================================
# Step 2: set values / formulas for the new cells (rows 100-103)
# ============================================================================
$ws1.Range("A100").Value = 46040
$ws1.Range("C100").Value = 4
$ws1.Range("D100").Value = 206603
$ws1.Range("E100").Value = 196
$ws1.Range("F100").Formula = "=D100/E100*1000"
$ws1.Range("G100").Formula = "=(E96-E100)/E96"
$ws1.Range("H100").Formula = "=(F100-80000000)/80000000"
$ws1.Range("I100").Value = 4
$ws1.Range("J100").Value = 197281
$ws1.Range("K100").Value = 7
$ws1.Range("L100").Formula = "=J100/K100*1000"
$ws1.Range("M100").Formula = "=(K96-K100)/K96"
$ws1.Range("N100").Formula = "=(L100-80000000)/80000000"
$ws1.Range("P100").Value = "Replaced cache with transposition table with Zobrist key"

$ws1.Range("C101").Value = 5
$ws1.Range("D101").Value = 5072212
$ws1.Range("E101").Value = 4121
$ws1.Range("F101").Formula = "=D101/E101*1000"
$ws1.Range("G101").Formula = "=(E97-E101)/E97"
$ws1.Range("H101").Formula = "=(F101-80000000)/80000000"
$ws1.Range("I101").Value = 5
$ws1.Range("J101").Value = 4880523
$ws1.Range("K101").Value = 145
$ws1.Range("L101").Formula = "=J101/K101*1000"
$ws1.Range("M101").Formula = "=(K97-K101)/K97"
$ws1.Range("N101").Formula = "=(L101-80000000)/80000000"

$ws1.Range("I102").Value = 6
$ws1.Range("J102").Value = 119060324
$ws1.Range("K102").Value = 3136
$ws1.Range("L102").Formula = "=J102/K102*1000"
$ws1.Range("M102").Formula = "=(K98-K102)/K98"
$ws1.Range("N102").Formula = "=(L102-80000000)/80000000"

$ws1.Range("I103").Value = 7
$ws1.Range("J103").Value = 3195901860
$ws1.Range("K103").Value = 69384
$ws1.Range("L103").Formula = "=J103/K103*1000"
$ws1.Range("N103").Formula = "=(L103-80000000)/80000000"

# ============================================================================
# Step 3: re-style the now-superseded "latest" markers P96 / P97 from the
# highlighted note style down to the plain note style (s=38 -> s=7)
# ============================================================================
$ws1.Range("P6").Copy()
$ws1.Range("P96").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("P6").Copy()
$ws1.Range("P97").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ============================================================================
# Step 4: update the view state - active cell moves to the newest entry
# ============================================================================
$ws1.Activate()
$ws1.Range("A56").Select()
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("L103").Select()
